$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "SCRIPT/D24P11A/enter06.ssb"
$ws.Range("A12").Value = "SCRIPT/D25P11A/enter06.ssb"
$ws.Range("A13").Value = "SCRIPT/D27P11A/enter02.ssb"

$ws.Range("A11:A13").WrapText = $true

$ws.Range("B4").Select()
